$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

$row2 = @(0.049644077596912715, 0.80710424214941401, 8.173, 143.76599999999999, 33.56, 4.7539999999999996, 5.0810000000000004, 0.048142803698359257, 0.80665040903704677, 43.765999999999998, 143.76599999999999, 4.5339999999999998, 4.851)
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, 3 + $i).Value = $row2[$i] }
$row3 = @(0.052195579972023488, 0.8100170689821985, 8.1969999999999992, 143.93600000000001, 33.652000000000001, 3.4870000000000001, 3.452, 0.056201550387596902, 0.80917898319413428, 43.936, 104.447, 3.3959999999999999, 3.38)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, 3 + $i).Value = $row3[$i] }
$row4 = @(0.049208736139267316, 0.80554951202278147, 8.1920000000000002, 143.67699999999999, 33.533999999999999, 4.9329999999999998, 4.5960000000000001, 0.037527593818984552, 0.80624788124788116, 43.677, 106.18, 4.7229999999999999, 4.7080000000000002)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, 3 + $i).Value = $row4[$i] }
$row5 = @(0.050100961200887972, 0.80586212794299705, 8.1950000000000003, 143.71100000000001, 33.399000000000001, 7.9279999999999999, 8, 0.035219399538106239, 0.80609014183271088, 43.710999999999999, 110.387, 7.94, 7.9180000000000001)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, 3 + $i).Value = $row5[$i] }
$row6 = @(0.076717898371446755, 0.88523034434975523, 8.1809999999999992, 149.05099999999999, 40.323999999999998, 5.9169999999999998, 5.6139999999999999, 0.074273715378009933, 0.8860799904169081, 49.051000000000002, 149.05099999999999, 5.9379999999999997, 5.42)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, 3 + $i).Value = $row6[$i] }
$row7 = @(0.077183039448961904, 0.88637864517715204, 8.1709999999999994, 149.02000000000001, 40.32, 4.0979999999999999, 3.9689999999999999, 0.090507726269315664, 0.88576165060068435, 49.02, 105.03100000000001, 4.1550000000000002, 4.1230000000000002)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, 3 + $i).Value = $row7[$i] }
$row8 = @(0.075718584904645475, 0.88664906695775769, 8.1890000000000001, 149.089, 40.445, 5.98, 6.3390000000000004, 0.093648867313915862, 0.88194729861396537, 49.088999999999999, 106.979, 5.7030000000000003, 5.6539999999999999)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, 3 + $i).Value = $row8[$i] }
$row9 = @(0.07471772862879357, 0.88498917362520524, 8.1750000000000007, 148.78100000000001, 40.231999999999999, 8.9689999999999994, 8.452, 0.081771720613287913, 0.87880869595231037, 48.780999999999999, 111.538, 9.4239999999999995, 9.3610000000000007)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, 3 + $i).Value = $row9[$i] }
$row10 = @(0.10243691593213265, 0.96748770967190523, 8.1820000000000004, 154.80799999999999, 48.816000000000003, 6.9450000000000003, 7.3120000000000003, 0.097909891290031817, 0.96895602456638286, 54.808, 154.80799999999999, 6.819, 6.6559999999999997)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, 3 + $i).Value = $row10[$i] }
$row11 = @(0.10466067803700209, 0.96901931149829157, 8.1850000000000005, 154.63399999999999, 48.48, 4.9589999999999996, 4.9930000000000003, 0.079944178628389165, 0.97159718694322506, 54.634, 105.565, 4.9349999999999996, 4.8860000000000001)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, 3 + $i).Value = $row11[$i] }
$row12 = @(0.10297371591668587, 0.96833594794039302, 8.1769999999999996, 154.57599999999999, 48.576000000000001, 6.7210000000000001, 6.4960000000000004, 0.094986168741355459, 0.96934135586635595, 54.576000000000001, 107.876, 7.0309999999999997, 6.9660000000000002)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws.Cells.Item(12, 3 + $i).Value = $row12[$i] }
$row13 = @(0.10676687178942722, 0.96926261352211363, 8.1780000000000008, 154.905, 48.692999999999998, 12.515000000000001, 11.177, 0.10818965517241383, 0.96973036460592887, 54.905000000000001, 112.95699999999999, 11.535, 11.404999999999999)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws.Cells.Item(13, 3 + $i).Value = $row13[$i] }
$row14 = @(0.11892534532002899, 0.99100391450142367, 8.1880000000000006, 156.404, 50.741999999999997, 7.202, 7.4550000000000001, 0.1216608105496997, 0.98983400626004692, 56.404000000000003, 156.404, 7.3410000000000002, 7.1790000000000003)
for ($i = 0; $i -lt $row14.Length; $i++) { $ws.Cells.Item(14, 3 + $i).Value = $row14[$i] }
$row15 = @(0.11752976547397881, 0.9910788597548208, 8.1940000000000008, 156.548, 51.012999999999998, 5.1210000000000004, 5.1109999999999998, 0.10736800630417651, 0.98895392557364359, 56.548000000000002, 105.681, 5.1360000000000001, 5.07)
for ($i = 0; $i -lt $row15.Length; $i++) { $ws.Cells.Item(15, 3 + $i).Value = $row15[$i] }
$row16 = @(0.11995540633948878, 0.99096104043937716, 8.2010000000000005, 156.749, 51.207000000000001, 7.3620000000000001, 7.1379999999999999, 0.11528973034997134, 0.98985498112998116, 56.749000000000002, 108.31100000000001, 7.4770000000000003, 7.3739999999999997)
for ($i = 0; $i -lt $row16.Length; $i++) { $ws.Cells.Item(16, 3 + $i).Value = $row16[$i] }
$row17 = @(0.11872876605133117, 0.99019341955146667, 8.1760000000000002, 156.26300000000001, 50.698999999999998, 12.134, 11.17, 0.13212401055408968, 0.99026358088505861, 56.262999999999998, 113.334, 11.987, 11.816000000000001)
for ($i = 0; $i -lt $row17.Length; $i++) { $ws.Cells.Item(17, 3 + $i).Value = $row17[$i] }

Write-Output "done"